$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure cells keep their original text representation (e.g. "63.60", "1.00")
# instead of being auto-converted to numbers by Excel when the value looks numeric.
$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "26.182.99"
$ws.Range("E2").NumberFormat = "@"
$ws.Range("E2").Value = "  -0.51%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.591.25"
$ws.Range("E4").NumberFormat = "@"
$ws.Range("E4").Value = "  -0.08%  "
$ws.Range("E5").NumberFormat = "@"
$ws.Range("E5").Value = "  +0.84%  "
$ws.Range("E6").NumberFormat = "@"
$ws.Range("E6").Value = "  -0.88%  "
$ws.Range("E8").NumberFormat = "@"
$ws.Range("E8").Value = "  -0.12%  "
$ws.Range("E9").NumberFormat = "@"
$ws.Range("E9").Value = "  -0.88%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "18.98"
$ws.Range("E10").NumberFormat = "@"
$ws.Range("E10").Value = "  -2.17%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0845"
$ws.Range("E11").NumberFormat = "@"
$ws.Range("E11").Value = "  +0.05%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "1.814.90"
$ws.Range("E12").NumberFormat = "@"
$ws.Range("E12").Value = "  +0.16%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "1.606.97"
$ws.Range("E13").NumberFormat = "@"
$ws.Range("E13").Value = "  +1.17%  "
$ws.Range("E14").NumberFormat = "@"
$ws.Range("E14").Value = "  -1.35%  "
$ws.Range("E15").NumberFormat = "@"
$ws.Range("E15").Value = "  -1.97%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "63.60"
$ws.Range("E16").NumberFormat = "@"
$ws.Range("E16").Value = "  -1.18%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "26.190.17"
$ws.Range("E17").NumberFormat = "@"
$ws.Range("E17").Value = "  -0.50%  "
$ws.Range("E18").NumberFormat = "@"
$ws.Range("E18").Value = "  -0.50%  "
$ws.Range("B19").NumberFormat = "@"
$ws.Range("B19").Value = "Chainlink"
$ws.Range("C19").NumberFormat = "@"
$ws.Range("C19").Value = "https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link"
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.34"
$ws.Range("E19").NumberFormat = "@"
$ws.Range("E19").Value = "  -1.48%  "
$ws.Range("B20").NumberFormat = "@"
$ws.Range("B20").Value = "BitcoinCash"
$ws.Range("C20").NumberFormat = "@"
$ws.Range("C20").Value = "https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch"
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "213.61"
$ws.Range("E20").NumberFormat = "@"
$ws.Range("E20").Value = "  +1.26%  "
$ws.Range("E21").NumberFormat = "@"
$ws.Range("E21").Value = "  +0.01%  "
$ws.Range("E22").NumberFormat = "@"
$ws.Range("E22").Value = "  -0.63%  "
$ws.Range("E23").NumberFormat = "@"
$ws.Range("E23").Value = "  +0.84%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "144.81"
$ws.Range("E25").NumberFormat = "@"
$ws.Range("E25").Value = "  +0.20%  "
$ws.Range("E26").NumberFormat = "@"
$ws.Range("E26").Value = "  -0.08%  "
$ws.Range("E28").NumberFormat = "@"
$ws.Range("E28").Value = "  -1.28%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "15.08"
$ws.Range("E29").NumberFormat = "@"
$ws.Range("E29").Value = "  -1.10%  "
$ws.Range("E30").NumberFormat = "@"
$ws.Range("E30").Value = "  -2.56%  "
$ws.Range("E31").NumberFormat = "@"
$ws.Range("E31").Value = "  +0.40%  "
$ws.Range("E32").NumberFormat = "@"
$ws.Range("E32").Value = "  -1.68%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.421.61"
$ws.Range("E33").NumberFormat = "@"
$ws.Range("E33").Value = "  +7.94%  "
$ws.Range("E34").NumberFormat = "@"
$ws.Range("E34").Value = "  -1.65%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.43"
$ws.Range("E35").NumberFormat = "@"
$ws.Range("E35").Value = "  -0.78%  "
$ws.Range("E36").NumberFormat = "@"
$ws.Range("E36").Value = "  -0.85%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.586"
$ws.Range("E37").NumberFormat = "@"
$ws.Range("E37").Value = "  -3.84%  "
$ws.Range("E38").NumberFormat = "@"
$ws.Range("E38").Value = "  -1.75%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.91"
$ws.Range("E39").NumberFormat = "@"
$ws.Range("E39").Value = "  +4.83%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.822"
$ws.Range("E40").NumberFormat = "@"
$ws.Range("E40").Value = "  +2.15%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "1.00"
$ws.Range("E41").NumberFormat = "@"
$ws.Range("E41").Value = "  -0.05%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.974"
$ws.Range("E42").NumberFormat = "@"
$ws.Range("E42").Value = "  -9.85%  "
$ws.Range("E43").NumberFormat = "@"
$ws.Range("E43").Value = "  -0.21%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "2.12"
$ws.Range("E44").NumberFormat = "@"
$ws.Range("E44").Value = "  -0.35%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "1.726.14"
$ws.Range("E45").NumberFormat = "@"
$ws.Range("E45").Value = "  +0.06%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "60.94"
$ws.Range("E46").NumberFormat = "@"
$ws.Range("E46").Value = "  -2.19%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "86.95"
$ws.Range("E47").NumberFormat = "@"
$ws.Range("E47").Value = "  -0.65%  "
$ws.Range("B48").NumberFormat = "@"
$ws.Range("B48").Value = "RenderToken"
$ws.Range("C48").NumberFormat = "@"
$ws.Range("C48").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.48"
$ws.Range("E48").NumberFormat = "@"
$ws.Range("E48").Value = "  -0.27%  "
$ws.Range("B49").NumberFormat = "@"
$ws.Range("B49").Value = "Cronos"
$ws.Range("C49").NumberFormat = "@"
$ws.Range("C49").Value = "https://coinranking.com/coin/65PHZTpmE55b+cronos-cro"
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.0502"
$ws.Range("E49").NumberFormat = "@"
$ws.Range("E49").Value = "  -0.91%  "
$ws.Range("B50").NumberFormat = "@"
$ws.Range("B50").Value = "Algorand"
$ws.Range("C50").NumberFormat = "@"
$ws.Range("C50").Value = "https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo"
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "0.0957"
$ws.Range("E50").NumberFormat = "@"
$ws.Range("E50").Value = "  -2.01%  "
$ws.Range("B51").NumberFormat = "@"
$ws.Range("B51").Value = "USDD"
$ws.Range("C51").NumberFormat = "@"
$ws.Range("C51").Value = "https://coinranking.com/coin/z2PZIKQL7+usdd-usdd"
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.998"
$ws.Range("E51").NumberFormat = "@"
$ws.Range("E51").Value = "  -0.20%  "
